$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1, matching the formatting of the other header cells (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add the time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 10:52:42.963362"
$ws.Range("F3").Value = "2021-10-05 10:52:42.963375"
$ws.Range("F4").Value = "2021-10-05 10:52:42.963379"
